$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A values (rows 2-11) from 1 to 3
$ws.Range("A2:A11").Value = 3

# Update the selection to K16 (single cell, no range)
$ws.Range("K16").Select()
